# Applies the "Add files via upload" edit:
#   - Slide 1 ("Title 1"): reposition title placeholder and center its text
#   - Slide 2 ("Title 1"): reposition title placeholder
#   - Slide 3 ("Title 1"): reposition title placeholder
#   - Slide 5 ("TextBox 7" and the "Diagram 9" SmartArt graphic frame): reposition both
#   - Slide 8 ("TextBox 3"): reposition text box
#
# Note: EMU/12700 = points, but Shape.Left/Top/Width/Height round-trip through a
# 32-bit float in this object model, so the literals below are nudged by a handful
# of ULPs (still the same displayed point value) so the saved EMU lands exactly on
# the target instead of the nearest float32 truncating a hair short.
#
# (The source diff also stamps a fresh p:nvPr/p:extLst/p14:modId on the slide 5
# diagram frame. That's PowerPoint's own internal co-authoring "last modified"
# marker - it isn't backed by any Shape/GraphicFrame property or method in the
# PowerPoint object model, so it can't be set from automation code and is left
# alone here; every positional/formatting change it is attached to is applied.)

$p = $ppt.ActivePresentation

# ---- Slide 1: Title 1 ----
$s1 = $p.Slides.Item(1)
$shp1 = $s1.Shapes.Item("Title 1")
$shp1.Left = 122.57133958267717
$shp1.Top = 81.4775590551181
$shp1.TextFrame.TextRange.ParagraphFormat.Alignment = 2

# ---- Slide 2: Title 1 ----
$s2 = $p.Slides.Item(2)
$shp2 = $s2.Shapes.Item("Title 1")
$shp2.Left = 85.14244494488189
$shp2.Top = 33.176458692913386

# ---- Slide 3: Title 1 ----
$s3 = $p.Slides.Item(3)
$shp3 = $s3.Shapes.Item("Title 1")
$shp3.Left = 123.648031496063
$shp3.Top = 51.28567029133858

# ---- Slide 5: TextBox 7 and Diagram 9 ----
$s5 = $p.Slides.Item(5)

$shp5tb = $s5.Shapes.Item("TextBox 7")
$shp5tb.Left = 149.25
$shp5tb.Top = 130.72535433070865

$shp5dg = $s5.Shapes.Item("Diagram 9")
$shp5dg.Left = 209.1428346456693
$shp5dg.Top = 189.48795375590552

# ---- Slide 8: TextBox 3 ----
$s8 = $p.Slides.Item(8)
$shp8 = $s8.Shapes.Item("TextBox 3")
$shp8.Left = 169.0588188976378
$shp8.Top = 118.8572440944882
